$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellD = $ws.Cells.Item(2, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "25.478.35"
$cellD.ClearFormats()
$ws.Range("E2").Value = "  +1.56%  "

$cellD = $ws.Cells.Item(3, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.661.55"
$cellD.ClearFormats()
$ws.Range("E3").Value = "  +0.56%  "

$cellD = $ws.Cells.Item(4, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.9995"
$cellD.ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "

$cellD = $ws.Cells.Item(5, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "236.96"
$cellD.ClearFormats()
$ws.Range("E5").Value = "  -0.71%  "

$ws.Range("E6").Value = "  +0.01%  "

$cellD = $ws.Cells.Item(7, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.4785"
$cellD.ClearFormats()
$ws.Range("E7").Value = "  +0.25%  "

$cellD = $ws.Cells.Item(8, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.2614"
$cellD.ClearFormats()
$ws.Range("E8").Value = "  +0.25%  "

$cellD = $ws.Cells.Item(9, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.06164"
$cellD.ClearFormats()
$ws.Range("E9").Value = "  +2.63%  "

$cellD = $ws.Cells.Item(10, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.07101"
$cellD.ClearFormats()
$ws.Range("E10").Value = "  -1.00%  "

$cellD = $ws.Cells.Item(11, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.659.38"
$cellD.ClearFormats()
$ws.Range("E11").Value = "  +0.47%  "

$cellD = $ws.Cells.Item(12, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "14.74"
$cellD.ClearFormats()
$ws.Range("E12").Value = "  +1.14%  "

$cellD = $ws.Cells.Item(13, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.5868"
$cellD.ClearFormats()
$ws.Range("E13").Value = "  -5.56%  "

$cellD = $ws.Cells.Item(14, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "4.356"
$cellD.ClearFormats()
$ws.Range("E14").Value = "  -5.11%  "

$cellD = $ws.Cells.Item(15, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "74.48"
$cellD.ClearFormats()
$ws.Range("E15").Value = "  +1.64%  "

$ws.Range("E16").Value = "  -0.06%  "

$cellD = $ws.Cells.Item(17, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.000"
$cellD.ClearFormats()
$ws.Range("E17").Value = "  -0.01%  "

$cellD = $ws.Cells.Item(18, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "25.482.24"
$cellD.ClearFormats()
$ws.Range("E18").Value = "  +1.64%  "

$cellD = $ws.Cells.Item(19, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.000006757"
$cellD.ClearFormats()
$ws.Range("E19").Value = "  +2.69%  "

$cellD = $ws.Cells.Item(20, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "11.42"
$cellD.ClearFormats()
$ws.Range("E20").Value = "  +0.25%  "

$cellD = $ws.Cells.Item(21, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.868.80"
$cellD.ClearFormats()
$ws.Range("E21").Value = "  +0.59%  "

$cellD = $ws.Cells.Item(22, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "4.431"
$cellD.ClearFormats()
$ws.Range("E22").Value = "  +0.07%  "

$cellD = $ws.Cells.Item(23, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "8.646"
$cellD.ClearFormats()
$ws.Range("E23").Value = "  +0.37%  "

$cellD = $ws.Cells.Item(24, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "5.256"
$cellD.ClearFormats()
$ws.Range("E24").Value = "  -0.16%  "

$cellD = $ws.Cells.Item(25, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "132.76"
$cellD.ClearFormats()
$ws.Range("E25").Value = "  +0.14%  "

$cellD = $ws.Cells.Item(26, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "15.00"
$cellD.ClearFormats()
$ws.Range("E26").Value = "  +0.98%  "

$cellD = $ws.Cells.Item(27, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.386"
$cellD.ClearFormats()
$ws.Range("E27").Value = "  -0.14%  "

$cellD = $ws.Cells.Item(28, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "104.74"
$cellD.ClearFormats()
$ws.Range("E28").Value = "  +1.89%  "

$cellD = $ws.Cells.Item(29, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.696"
$cellD.ClearFormats()
$ws.Range("E29").Value = "  +1.43%  "

$cellD = $ws.Cells.Item(30, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "3.927"
$cellD.ClearFormats()
$ws.Range("E30").Value = "  +4.22%  "

$cellD = $ws.Cells.Item(31, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "3.647"
$cellD.ClearFormats()
$ws.Range("E31").Value = "  +1.50%  "

$cellD = $ws.Cells.Item(32, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.07604"
$cellD.ClearFormats()
$ws.Range("E32").Value = "  -3.69%  "

$cellD = $ws.Cells.Item(33, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.9994"
$cellD.ClearFormats()
$ws.Range("E33").Value = "  +0.00%  "

$cellD = $ws.Cells.Item(34, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.04219"
$cellD.ClearFormats()
$ws.Range("E34").Value = "  -7.94%  "

$cellD = $ws.Cells.Item(35, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "2.615"
$cellD.ClearFormats()
$ws.Range("E35").Value = "  +0.26%  "

$cellD = $ws.Cells.Item(36, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.6089"
$cellD.ClearFormats()
$ws.Range("E36").Value = "  +6.24%  "

$cellD = $ws.Cells.Item(37, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.9501"
$cellD.ClearFormats()
$ws.Range("E37").Value = "  +1.12%  "

$cellD = $ws.Cells.Item(38, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "2.602"
$cellD.ClearFormats()
$ws.Range("E38").Value = "  -0.43%  "

$cellD = $ws.Cells.Item(39, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.8606"
$cellD.ClearFormats()
$ws.Range("E39").Value = "  +1.83%  "

$cellD = $ws.Cells.Item(40, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.9996"
$cellD.ClearFormats()
$ws.Range("E40").Value = "  -0.07%  "

$cellD = $ws.Cells.Item(41, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.850"
$cellD.ClearFormats()
$ws.Range("E41").Value = "  +1.06%  "

$cellD = $ws.Cells.Item(42, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.01469"
$cellD.ClearFormats()
$ws.Range("E42").Value = "  -5.00%  "

$cellD = $ws.Cells.Item(43, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "97.04"
$cellD.ClearFormats()
$ws.Range("E43").Value = "  -1.68%  "

$cellD = $ws.Cells.Item(44, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.3748"
$cellD.ClearFormats()
$ws.Range("E44").Value = "  +1.53%  "

$cellD = $ws.Cells.Item(45, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "4.723"
$cellD.ClearFormats()
$ws.Range("E45").Value = "  -1.43%  "

$cellD = $ws.Cells.Item(46, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.1117"
$cellD.ClearFormats()
$ws.Range("E46").Value = "  -0.90%  "

$cellD = $ws.Cells.Item(47, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "6.198"
$cellD.ClearFormats()
$ws.Range("E47").Value = "  +2.18%  "

$cellD = $ws.Cells.Item(48, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "0.05251"
$cellD.ClearFormats()
$ws.Range("E48").Value = "  +1.32%  "

$cellD = $ws.Cells.Item(49, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "29.51"
$cellD.ClearFormats()
$ws.Range("E49").Value = "  -0.49%  "

$ws.Range("E50").Value = "  +0.03%  "

$cellD = $ws.Cells.Item(51, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "1.001"
$cellD.ClearFormats()
$ws.Range("E51").Value = "  +0.04%  "
